$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are free-form text (prices with thousands separators,
# trailing zeros, subscript digits, etc). Some of them look like numbers
# (e.g. '560.80', '5.30', '0.999'), so a leading apostrophe is used to force
# Excel to keep them as text, and the style is reset afterwards so no
# quote-prefix formatting lingers on the cell.

$ws.Range('D2').Value = "'64.078.56"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = "'3.148.69"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'560.80"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'140.69"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = "'3.142.88"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.97%  '
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('D10').Value = "'6.70"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').Value = "'36.21"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').Value = "'3.660.34"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.26%  '
$ws.Range('D16').Value = "'64.199.62"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').Value = "'3.151.29"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.94%  '
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').Value = "'509.53"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.31%  '
$ws.Range('D20').Value = "'6.77"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('D21').Value = "'13.94"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('E22').Value = '  +2.73%  '
$ws.Range('D23').Value = "'7.39"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('E24').Value = '  +1.57%  '
$ws.Range('D25').Value = "'78.47"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').Value = "'8.68"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.84%  '
$ws.Range('E28').Value = '  +2.49%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').Value = "'0.999"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').Value = "'26.50"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').Value = "'551.05"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.47%  '
$ws.Range('D35').Value = "'6.04"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = "'53.77"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.05%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = "'5.30"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.74%  '
$ws.Range('D38').Value = "'0.0424"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.70%  '
$ws.Range('D39').Value = "'3.152.99"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.19%  '
$ws.Range('D40').Value = "'0.0815"
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'0.121"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('E42').Value = '  -7.21%  '
$ws.Range('D43').Value = "'8.22"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.70%  '
$ws.Range('E44').Value = '  +6.27%  '
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('D47').Value = "'121.87"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.05%  '
$ws.Range('D48').Value = "'24.77"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.56%  '
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').Value = "'0.0₃0511"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.47%  '
$ws.Range('E51').Value = '  -1.05%  '
